# "append desc to days in week"
# - Template sheet: scroll/selection view-state tweak (tabSelected -> topLeftCell area),
#   selection stays at B10.
# - May sheet: remove the stray "It Worked!!" note (J2), rewrite the week rows with
#   per-day description rows (week1..week4 headers followed by numbered day rows),
#   and make "May" the active/selected sheet with the selection on I7.

$wb = $excel.ActiveWorkbook

$template = $wb.Worksheets.Item("Template")
$may = $wb.Worksheets.Item("May")

# --- Template sheet: keep selection on B10, just touch the view (no data changes) ---
$template.Activate()
$template.Range("B10").Select() | Out-Null

# --- May sheet: drop the leftover "It Worked!!" cell ---
$may.Range("J2").ClearContents() | Out-Null

# Row 7 becomes a one-off scratch/day row (no longer "week 1" header text)
$may.Range("A7").Value = "22"
$may.Range("B7").Value = "efwef"
$may.Range("C7").Value = "*"

# Week headers (column A only) followed by 5 numbered day rows each (B = day index, C = "*")
$weekHeaderRows = @(8, 14, 20, 26)
$weekNames = @("week1", "week2", "week3", "week4")

for ($w = 0; $w -lt 4; $w++) {
    $headerRow = $weekHeaderRows[$w]
    $may.Range("A$headerRow").Value = $weekNames[$w]
    $may.Range("B$headerRow`:C$headerRow").ClearContents() | Out-Null

    $dayStart = $headerRow + 1
    for ($i = 0; $i -lt 5; $i++) {
        $r = $dayStart + $i
        $dayNumber = ($w * 5) + $i + 1
        $may.Range("A$r").Value = "22"
        $may.Range("B$r").Value = [string]$dayNumber
        $may.Range("C$r").Value = "*"
    }
}

# Make "May" the active sheet with the new selection
$may.Activate()
$may.Range("I7").Select() | Out-Null
